$wb = $excel.ActiveWorkbook

# Remove the second worksheet entirely (ID_9a601b4 row is dropped)
$wb.Worksheets.Item(2).Delete()

# Rename the remaining worksheet
$ws = $wb.Worksheets.Item(1)
$ws.Name = "ID_dda333c"

# Update the data row with the new product details
$ws.Range("A2").Value = 22.3
$ws.Range("B2").Value = "23/02/2025"
$ws.Range("C2").Value = "Mens Casual Premium Slim Fit T-Shirts "
$ws.Range("D2").Value = "Slim-fitting style, contrast raglan long sleeve, three-button henley placket, light weight & soft fabric for breathable and comfortable wearing. And Solid stitched shirts with round neck made for durability and a great fit for casual fashion wear and diehard baseball fans. The Henley style round neckline includes a three-button placket."
$ws.Range("E2").Value = "ID_dda333c"
